$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(1)
$textRange = $shape.TextFrame.TextRange

# The paragraph currently holds three separate runs: "Below", " " and
# "section-level". Re-assigning the identical concatenated text is a
# no-op for the runtime, so first nudge the text to a different value
# and then set the desired final text; this forces the runs to be
# rebuilt (and consolidated) as a single run.
$textRange.Text = "-"
$textRange.Text = "Below section-level"
